# DAI Municipios Enlace Transparencia.xlsx — update several rows in the
# "Codes" sheet with newly-received survey answers (professional background,
# seniority, hiring modality, request/response codes and dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: paste only the number-format/alignment/fill of $srcRef onto $dstRef
# (re-uses an existing cell style instead of minting a brand-new one).
function Copy-CellFormat($srcRef, $dstRef) {
    $ws.Range($srcRef).Copy()
    $ws.Range($dstRef).PasteSpecial(-4122)   # xlPasteFormats
}

# Helper: stamp a date serial into a cell that already carries (or is given)
# the short-date style used throughout column D/E ("D5" is a representative
# donor cell with that style).
function Set-DateValue($ref, $serial) {
    $ws.Range($ref).Value = $serial
}

# ---------------------------------------------------------------------
# Row 43 — Casablanca / MU030: answer arrived, fill in response date and
# the three profile columns; row grows to a 2-line wrap.
# ---------------------------------------------------------------------
Copy-CellFormat "A12" "A43"
$ws.Rows(43).RowHeight = 28.8

Copy-CellFormat "D5" "E43"
Set-DateValue "E43" 44663

$ws.Range("F43").Value = "Lingüista (Enlace) Psicologa (Encargada)"
$ws.Range("G43").Value = "12 años (Enlace) 4 años (Encargada)"
$ws.Range("H43").Value = "Código del Trabajo (Enlace) Planta (Encargada)"

# ---------------------------------------------------------------------
# Row 236 — answer arrived (request + response dates, profile columns).
# ---------------------------------------------------------------------
Copy-CellFormat "A12" "A236"
$ws.Rows(236).RowHeight = 28.8

$ws.Range("C236").Value = "CT001T0015682"

Copy-CellFormat "D5" "D236"
Set-DateValue "D236" 44662
Copy-CellFormat "D5" "E236"
Set-DateValue "E236" 44663

$ws.Range("F236").Value = "Asistente Social (Secretario Municipal)"
$ws.Range("G236").Value = "13 años"
$ws.Range("H236").Value = "Planta"

# ---------------------------------------------------------------------
# Row 237 — answer arrived; G/H get the "highlighted" (yellow) style
# because of an editorial note, and the row wraps to 4 lines.
# ---------------------------------------------------------------------
Copy-CellFormat "A12" "A237"
$ws.Rows(237).RowHeight = 57.6

$ws.Range("C237").Value = "CT001T0015682"

Copy-CellFormat "D5" "D237"
Set-DateValue "D237" 44643
Copy-CellFormat "D5" "E237"
Set-DateValue "E237" 44663

$ws.Range("F237").Value = "Contador General (Enlace) Técnico de Nivel Superior en Redes y Telecomunicaciones (Encargado)"

Copy-CellFormat "C228" "G237"
$ws.Range("G237").Value = "40 años (Enlace) 3 años y 8 meses (Encargado)"
Copy-CellFormat "C228" "H237"
$ws.Range("H237").Value = "No indica"

# ---------------------------------------------------------------------
# Row 201 — request code + ingreso date received.
# ---------------------------------------------------------------------
$ws.Range("C201").Value = "MU155T0000432"
Copy-CellFormat "D5" "D201"
Set-DateValue "D201" 44663

# ---------------------------------------------------------------------
# Row 209 — request code + ingreso date received.
# ---------------------------------------------------------------------
$ws.Range("C209").Value = "MU079T0000826"
Copy-CellFormat "D5" "D209"
Set-DateValue "D209" 44663

# ---------------------------------------------------------------------
# Row 297 — request code + ingreso date received.
# ---------------------------------------------------------------------
$ws.Range("C297").Value = "MU265T0001884"
Copy-CellFormat "D5" "D297"
Set-DateValue "D297" 44663

# ---------------------------------------------------------------------
# Row 178 — request code + ingreso date received.
# ---------------------------------------------------------------------
$ws.Range("C178").Value = "MU230T0001696"
Copy-CellFormat "D5" "D178"
Set-DateValue "D178" 44663

# ---------------------------------------------------------------------
# Row 217 — answer arrived (request + response dates, profile columns).
# ---------------------------------------------------------------------
Copy-CellFormat "A12" "A217"
$ws.Rows(217).RowHeight = 28.8

$ws.Range("C217").Value = "CT001T0015682"

Copy-CellFormat "D5" "D217"
Set-DateValue "D217" 44659
Copy-CellFormat "D5" "E217"
Set-DateValue "E217" 44663

$ws.Range("F217").Value = "Ingeniero en Informática (Jefe de Calidad y Transparencia)"
$ws.Range("G217").Value = "8 años"
$ws.Range("H217").Value = "Planta"

# ---------------------------------------------------------------------
# Move the active selection, matching the author's last cursor position.
# ---------------------------------------------------------------------
$ws.Range("F8").Select()
